# Renaming tabs for consistency
#  - "default_constants" -> "constants"
#  - move the recorded selection on "constants" from C5 to C50 (scrolled to A44)
#  - move the recorded selection on "dropdown_lists" from D3 to A2:D4 (active D4)

$wb = $excel.ActiveWorkbook

$wsConstants = $wb.Worksheets.Item("default_constants")
$wsConstants.Name = "constants"

$wsDropdown = $wb.Worksheets.Item("dropdown_lists")

# --- "constants" sheet: selection C5 -> C50, view scrolled so A44 is the top-left cell ---
$wsConstants.Activate() | Out-Null
$wsConstants.Range("C50").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 44
$excel.ActiveWindow.ScrollColumn = 1

# --- "dropdown_lists" sheet: selection D3 -> A2:D4 (active cell D4) ---
$wsDropdown.Activate() | Out-Null
$wsDropdown.Range("A2:D4").Select() | Out-Null

# Restore "constants" as the active/selected sheet (matches the original file)
$wsConstants.Activate() | Out-Null
$wsConstants.Range("C50").Select() | Out-Null

Write-Host "done"
